$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Find-ParagraphByText {
    param($doc, $text)
    $rng = $doc.Content
    $ok = $rng.Find.Execute($text, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        throw "Text not found: $text"
    }
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ($p.Range.Start -le $rng.Start -and $p.Range.End -ge $rng.End) {
            return $p
        }
    }
    throw "Paragraph not found for: $text"
}

# ---------------------------------------------------------------------------
# 1) After "Electrum Wallet" insert a new paragraph "Guide"
# ---------------------------------------------------------------------------
$pElectrum = Find-ParagraphByText $d "Electrum Wallet"
$pElectrum.Range.InsertParagraphAfter()
$pGuide = $pElectrum.Next()
$xmlGuide = @"
<w:p $wns>
  <w:pPr>
    <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
    <w:ind w:left="360" w:firstLine="360"/>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/>
      <w:color w:val="000000"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
      <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial" w:eastAsia="Times New Roman"/>
      <w:color w:val="000000"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
      <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
    </w:rPr>
    <w:tab/>
    <w:t>Guide</w:t>
  </w:r>
</w:p>
"@
$pGuide.Range.InsertXML($xmlGuide)

# ---------------------------------------------------------------------------
# 2) Update the "Graphics" paragraph's own formatting (Times New Roman pPr
#    mark -> Arial pPr mark), keep its run/text identical.
# ---------------------------------------------------------------------------
$pGraphics = Find-ParagraphByText $d "Graphics"
$xmlGraphics = @"
<w:p $wns>
  <w:pPr>
    <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
    <w:ind w:left="360" w:firstLine="360"/>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/>
      <w:color w:val="000000"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial" w:eastAsia="Times New Roman"/>
      <w:color w:val="000000"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>Graphics</w:t>
  </w:r>
</w:p>
"@
$pGraphics.Range.InsertXML($xmlGraphics)

# ---------------------------------------------------------------------------
# 2b) Insert "Marketing Materials" right after the "Graphics" paragraph
# ---------------------------------------------------------------------------
$pGraphics = Find-ParagraphByText $d "Graphics"
$pGraphics.Range.InsertParagraphAfter()
$pMarketing = $pGraphics.Next()
$xmlMarketing = @"
<w:p $wns>
  <w:pPr>
    <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
    <w:ind w:left="360" w:firstLine="360"/>
    <w:rPr>
      <w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman" w:eastAsia="Times New Roman"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:tab/>
    <w:t xml:space="preserve">Marketing Materials</w:t>
  </w:r>
</w:p>
"@
$pMarketing.Range.InsertXML($xmlMarketing)

# ---------------------------------------------------------------------------
# 3) Insert "Exchanges Listing Guide" right after the "Exchanges" paragraph
# ---------------------------------------------------------------------------
$pExchanges = Find-ParagraphByText $d "Exchanges"
$pExchanges.Range.InsertParagraphAfter()
$pExchGuide = $pExchanges.Next()
$xmlExchGuide = @"
<w:p $wns>
  <w:pPr>
    <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
    <w:ind w:left="360"/>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/>
      <w:color w:val="000000"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial" w:eastAsia="Times New Roman"/>
      <w:color w:val="000000"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:tab/>
    <w:tab/>
    <w:t xml:space="preserve">Exchanges Listing Guide</w:t>
  </w:r>
</w:p>
"@
$pExchGuide.Range.InsertXML($xmlExchGuide)

# ---------------------------------------------------------------------------
# 4) Insert "The Other Side" right after the "SmartCard" paragraph that
#    precedes "Resources" (the menu's first SmartCard occurrence).
# ---------------------------------------------------------------------------
$pResources = Find-ParagraphByText $d "Resources"
$pSmartCard = $pResources.Previous()
$pSmartCard.Range.InsertParagraphAfter()
$pOtherSide = $pSmartCard.Next()
$xmlOtherSide = @"
<w:p $wns>
  <w:pPr>
    <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
    <w:ind w:left="360"/>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/>
      <w:color w:val="000000"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial" w:eastAsia="Times New Roman"/>
      <w:color w:val="000000"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:tab/>
    <w:t xml:space="preserve">The Other Side</w:t>
  </w:r>
</w:p>
"@
$pOtherSide.Range.InsertXML($xmlOtherSide)

Write-Host "Done."
